# Revert responsive design implementation
# Re-adds the sensor data rows (A1:I28 -> A1:I31 on the FE-LIFTER sheets,
# and A1:I67 -> A1:I85 on the MID-LIFTER sheets) that a prior edit removed.

$wb = $excel.ActiveWorkbook

# Shared numeric value used by every new row's "ID_DEC" (col G) cell.
$gval = [double]"5.686312626471138e+23"

function Set-SensorRow($ws, $r, $timeVal, $bVal, $cVal, $dVal, $eVal, $fVal, $gVal, $hVal, $iVal) {
    $ws.Range("A$r").Value = $timeVal
    $ws.Range("A$r").NumberFormat = "YYYY-MM-DD HH:MM:SS"
    $ws.Range("B$r").Value = $bVal
    $ws.Range("C$r").Value = $cVal
    $ws.Range("D$r").Value = $dVal
    $ws.Range("E$r").Value = $eVal
    $ws.Range("F$r").Value = $fVal
    $ws.Range("G$r").Value = $gVal
    $ws.Range("H$r").Value = $hVal
    $ws.Range("I$r").Value = $iVal
}

# --- "FE-LIFTER" sheets: rows 29-31 (dimension A1:I28 -> A1:I31) ---
$feTimes = @(45729.58035023148, 45729.5803721875, 45729.58039546297)
$feRowStart = 29
$feB = "0x01,0x90"
$feC = "0x78,0x69,0x90,0x01,0x00,0x00,0x18,0x1b,0x41,0x0c,"
$feD = "0x01,0x90,"
$feE = "0x14"
$feF = 400
$feH = 400
$feI = 20

foreach ($sheetName in @("ROW50-FE-LIFTER", "ROW11-FE-LIFTER")) {
    $ws = $wb.Worksheets.Item($sheetName)
    for ($k = 0; $k -lt $feTimes.Count; $k++) {
        $r = $feRowStart + $k
        Set-SensorRow $ws $r $feTimes[$k] $feB $feC $feD $feE $feF $gval $feH $feI
    }
}

# --- "MID-LIFTER" sheets: rows 68-85 (dimension A1:I67 -> A1:I85) ---
$midTimes = @(
    45729.31518523148,
    45729.31520722222,
    45729.3152303588,
    45729.39866047454,
    45729.39868246527,
    45729.39870561343,
    45729.48213673611,
    45729.48215891204,
    45729.482181875,
    45729.56561206019,
    45729.56563403935,
    45729.56565724537,
    45729.64909206019,
    45729.64911011574,
    45729.64913337963,
    45729.73256368055,
    45729.73258552083,
    45729.73260887731
)
$midRowStart = 68
$midB = "0x01,0x90"
$midC = "0x78,0x69,0x90,0x01,0x00,0x00,0x18,0x20,0x41,0x0c,"
$midD = "0x01,0x90,"
$midE = "0x19"
$midF = 400
$midH = 400
$midI = 25

foreach ($sheetName in @("ROW50-MID-LIFTER", "ROW11-MID-LIFTER")) {
    $ws = $wb.Worksheets.Item($sheetName)
    for ($k = 0; $k -lt $midTimes.Count; $k++) {
        $r = $midRowStart + $k
        Set-SensorRow $ws $r $midTimes[$k] $midB $midC $midD $midE $midF $gval $midH $midI
    }
}
